# Auto-generated script applying targeted cell updates across multiple sheets
# of the Hades_Profits workbook, per the authoritative xml diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 207.91
$ws.Range("I15").Value = 207.91
$ws.Range("K15").Value = 623.73
$ws.Range("M15").Value = -454.73

$ws.Range("H86").Value = 1938.1333
$ws.Range("I86").Value = 1258
$ws.Range("J86").Value = 2715.4285
$ws.Range("K86").Value = 1258
$ws.Range("L86").Value = 2715.4285
$ws.Range("M86").Value = -135
$ws.Range("N86").Value = -4961.4285

$ws.Range("H89").Value = 1938.1333
$ws.Range("I89").Value = 1258
$ws.Range("J89").Value = 2715.4285
$ws.Range("K89").Value = 6290
$ws.Range("L89").Value = 13577.1425
$ws.Range("M89").Value = -674
$ws.Range("N89").Value = -24809.1425

$ws.Range("H131").Value = 1281.6097
$ws.Range("I131").Value = 889
$ws.Range("J131").Value = 1508.1154
$ws.Range("K131").Value = 2667
$ws.Range("L131").Value = 4524.3462
$ws.Range("M131").Value = 2373
$ws.Range("N131").Value = -14604.3462

$ws.Range("H140").Value = 53939.09
$ws.Range("J140").Value = 53939.09
$ws.Range("L140").Value = 53939.09
$ws.Range("N140").Value = -64299.09

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1630.02
$ws.Range("I32").Value = 1561.4082
$ws.Range("J32").Value = 4992
$ws.Range("K32").Value = 1561.4082
$ws.Range("L32").Value = 4992
$ws.Range("M32").Value = -1274.4082
$ws.Range("N32").Value = -5566

$ws.Range("H74").Value = 5082616.5
$ws.Range("I74").Value = 5765669
$ws.Range("J74").Value = 130487.5
$ws.Range("K74").Value = 5765669
$ws.Range("L74").Value = 130487.5
$ws.Range("M74").Value = -5764795
$ws.Range("N74").Value = -132235.5

$ws.Range("H77").Value = 5082616.5
$ws.Range("I77").Value = 5765669
$ws.Range("J77").Value = 130487.5
$ws.Range("K77").Value = 28828345
$ws.Range("L77").Value = 652437.5
$ws.Range("M77").Value = -28823977
$ws.Range("N77").Value = -661173.5

$ws.Range("H122").Value = 1273
$ws.Range("I122").Value = 1172.6111
$ws.Range("J122").Value = 1724.75
$ws.Range("K122").Value = 3517.8333
$ws.Range("L122").Value = 5174.25
$ws.Range("M122").Value = -1067.8333
$ws.Range("N122").Value = -10074.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 43890.445
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 43890.445
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 43890.445
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -44480.445

$ws.Range("H34").Value = 43890.445
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 43890.445
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 43890.445
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -44294.445

$ws.Range("H52").Value = 31766.666
$ws.Range("J52").Value = 32295.5
$ws.Range("L52").Value = 32295.5
$ws.Range("N52").Value = -32883.5

$ws.Range("H58").Value = 37038696
$ws.Range("I58").Value = 38463150
$ws.Range("J58").Value = 2800
$ws.Range("K58").Value = 38463150
$ws.Range("L58").Value = 2800
$ws.Range("M58").Value = -38462947
$ws.Range("N58").Value = -3206

$ws.Range("H132").Value = 45883.086
$ws.Range("I132").Value = 2292.0908
$ws.Range("J132").Value = 85841.5
$ws.Range("K132").Value = 6876.2724
$ws.Range("L132").Value = 257524.5
$ws.Range("M132").Value = -4346.2724
$ws.Range("N132").Value = -262584.5

$ws.Range("H136").Value = 37038696
$ws.Range("I136").Value = 38463150
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 115389450
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -115386900
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 527.94543
$ws.Range("J113").Value = 543.9375
$ws.Range("L113").Value = 1631.8125
$ws.Range("N113").Value = -5971.8125

$ws.Range("H131").Value = 895.13635
$ws.Range("I131").Value = 351.25
$ws.Range("J131").Value = 1016
$ws.Range("K131").Value = 1053.75
$ws.Range("L131").Value = 3048
$ws.Range("M131").Value = 3986.25
$ws.Range("N131").Value = -13128

$ws.Range("H133").Value = 6816.364
$ws.Range("J133").Value = 8050
$ws.Range("L133").Value = 24150
$ws.Range("N133").Value = -34270

$ws.Range("H140").Value = 2871.4634
$ws.Range("I140").Value = 3863.9473
$ws.Range("J140").Value = 2014.3182
$ws.Range("K140").Value = 11591.8419
$ws.Range("L140").Value = 6042.9546
$ws.Range("M140").Value = -6411.841899999999
$ws.Range("N140").Value = -16402.9546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.545456
$ws.Range("I2").Value = 34.4
$ws.Range("J2").Value = 141
$ws.Range("K2").Value = 34.4
$ws.Range("L2").Value = 141
$ws.Range("M2").Value = 78.59999999999999
$ws.Range("N2").Value = -367

$ws.Range("H80").Value = 4272.778
$ws.Range("I80").Value = 4010
$ws.Range("J80").Value = 4373.846
$ws.Range("K80").Value = 4010
$ws.Range("L80").Value = 4373.846
$ws.Range("M80").Value = -3012
$ws.Range("N80").Value = -6369.846

$ws.Range("H83").Value = 4272.778
$ws.Range("I83").Value = 4010
$ws.Range("J83").Value = 4373.846
$ws.Range("K83").Value = 20050
$ws.Range("L83").Value = 21869.23
$ws.Range("M83").Value = -15058
$ws.Range("N83").Value = -31853.23

$ws.Range("H132").Value = 39626.92
$ws.Range("I132").Value = 23117.822
$ws.Range("K132").Value = 69353.466
$ws.Range("M132").Value = -66823.466

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 724.63635
$ws.Range("I16").Value = 570.5263
$ws.Range("J16").Value = 1700.6666
$ws.Range("K16").Value = 570.5263
$ws.Range("L16").Value = 1700.6666
$ws.Range("M16").Value = -400.5263
$ws.Range("N16").Value = -2040.6666

$ws.Range("H46").Value = 2755496.2
$ws.Range("I46").Value = 3788488.5
$ws.Range("K46").Value = 3788488.5
$ws.Range("M46").Value = -3788300.5

$ws.Range("H61").Value = 2554.8333
$ws.Range("I61").Value = 2423.4546
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2423.4546
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2221.4546
$ws.Range("N61").Value = -4404

$ws.Range("H82").Value = 974.38464
$ws.Range("I82").Value = 733.36365
$ws.Range("K82").Value = 733.36365
$ws.Range("M82").Value = -372.36365

$ws.Range("H85").Value = 974.38464
$ws.Range("I85").Value = 733.36365
$ws.Range("K85").Value = 733.36365
$ws.Range("M85").Value = 514.63635

$ws.Range("H93").Value = 1114.9166
$ws.Range("I93").Value = 1029.6818
$ws.Range("J93").Value = 2052.5
$ws.Range("K93").Value = 1029.6818
$ws.Range("L93").Value = 2052.5
$ws.Range("M93").Value = 218.3181999999999
$ws.Range("N93").Value = -4548.5

$ws.Range("H100").Value = 1693.1177
$ws.Range("I100").Value = 1480.3
$ws.Range("K100").Value = 1480.3
$ws.Range("M100").Value = -939.3

$ws.Range("H113").Value = 2554.8333
$ws.Range("I113").Value = 2423.4546
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2423.4546
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -253.4546
$ws.Range("N113").Value = -8340

$ws.Range("H122").Value = 3472.5
$ws.Range("I122").Value = 3160
$ws.Range("J122").Value = 3597.5
$ws.Range("K122").Value = 9480
$ws.Range("L122").Value = 10792.5
$ws.Range("M122").Value = -7030
$ws.Range("N122").Value = -15692.5

$ws.Range("H136").Value = 70096.89999999999
$ws.Range("I136").Value = 40366.926
$ws.Range("J136").Value = 337666.66
$ws.Range("K136").Value = 121100.778
$ws.Range("L136").Value = 1012999.98
$ws.Range("M136").Value = -118550.778
$ws.Range("N136").Value = -1018099.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 440.21622
$ws.Range("I113").Value = 461.31818
$ws.Range("J113").Value = 409.26666
$ws.Range("K113").Value = 1383.95454
$ws.Range("L113").Value = 1227.79998
$ws.Range("M113").Value = 786.04546
$ws.Range("N113").Value = -5567.79998
